# Updated RT Individual Status Report
#
# The "Change Password" backlog item (row 64: # 18, "Change Password", 5,
# "JP, MD, RT, KW") moves down two rows to row 66, which was previously
# a blank spacer row between the "Iteration 2" and "Iteration 3" section
# headers. Row 64 becomes completely blank again afterwards.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Move the row-64 backlog entry down into row 66 (previously blank).
$ws.Range("A64:D64").Copy($ws.Range("A66:D66"))
$ws.Range("A64:D64").Clear()

# Match the reviewer's on-screen state when they saved: scrolled further
# down the sheet with the newly-relocated row selected as a whole row.
$ws.Activate()
$ws.Range("A66:XFD66").Select() | Out-Null
